$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-scrape data point (2026/01/26, 月, 19, 201) is inserted as a
# new row right before the existing row 719, shifting every row from the
# old 719 down by one (old 719 -> new 720, ..., old 760 -> new 761).
$ws.Rows.Item(719).Insert()

# Column A holds a date-like string ("2026/01/26") that must stay literal
# text (matching every other row in the sheet), not get auto-parsed into a
# date serial number. Force text formatting before the write, then restore
# the default "Normal" style so the cell matches its siblings exactly
# (no explicit style / number format).
$ws.Cells.Item(719, 1).NumberFormat = "@"
$ws.Cells.Item(719, 1).Value = "2026/01/26"
$ws.Cells.Item(719, 1).Style = "Normal"

$ws.Cells.Item(719, 2).Value = "月"
$ws.Cells.Item(719, 3).Value = 19
$ws.Cells.Item(719, 4).Value = 201
